$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.121.26"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "2.047.66"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'247.79"
$ws.Range("E5").Value = "  -2.47%  "
$ws.Range("E6").Value = "  -2.40%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'55.94"
$ws.Range("E8").Value = "  -7.30%  "
$ws.Range("D9").Value = "'0.379"
$ws.Range("E9").Value = "  -3.56%  "
$ws.Range("D10").Value = "'0.0778"
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("D11").Value = "'0.108"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "'16.22"
$ws.Range("D13").Value = "'0.878"
$ws.Range("E13").Value = "  +6.96%  "
$ws.Range("D14").Value = "2.345.61"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").Value = "'5.69"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "2.052.61"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "'18.37"
$ws.Range("E17").Value = "  +10.86%  "
$ws.Range("D18").Value = "37.111.19"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "'74.38"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  -4.66%  "
$ws.Range("E21").Value = "  -2.19%  "
$ws.Range("D22").Value = "'236.50"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").Value = "'9.51"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "'169.37"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").Value = "'2.16"
$ws.Range("E27").Value = "  -5.95%  "
$ws.Range("D28").Value = "'20.05"
$ws.Range("E28").Value = "  -2.18%  "
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "'4.85"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("D32").Value = "'0.0617"
$ws.Range("E32").Value = "  -3.65%  "
$ws.Range("D33").Value = "'4.47"
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("D34").Value = "'0.0886"
$ws.Range("E34").Value = "  -2.84%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  -3.70%  "
$ws.Range("D37").Value = "'1.78"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("D39").Value = "'5.25"
$ws.Range("E39").Value = "  +13.31%  "
$ws.Range("D40").Value = "'3.13"
$ws.Range("E40").Value = "  +9.40%  "
$ws.Range("D41").Value = "'0.0977"
$ws.Range("E41").Value = "  -17.74%  "
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'17.29"
$ws.Range("E43").Value = "  -3.98%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'1.15"
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("D45").Value = "'95.38"
$ws.Range("E45").Value = "  -4.18%  "
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").Value = "1.263.86"
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D49").Value = "'6.76"
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("D50").Value = "2.229.50"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").Value = "'43.75"
$ws.Range("E51").Value = "  -2.22%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
